$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(47).Insert()

$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(47, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value = "La Araucanía"
$ws.Cells.Item(47, 4).Value = 44476
$ws.Cells.Item(47, 5).Value = 9
$ws.Cells.Item(47, 6).Value = 100112035
$ws.Cells.Item(47, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 50
$ws.Cells.Item(47, 11).Value = 25000
$ws.Cells.Item(47, 12).Value = 25000
$ws.Cells.Item(47, 13).Value = 25000
$ws.Cells.Item(47, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(47, 16).Value = 2500
$ws.Cells.Item(47, 17).Value = 10
$ws.Cells.Item(47, 18).Value = "Hortaliza"
